$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $val)
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

$ws.Range('D2').Value = '38.006.20'
$ws.Range('E2').Value = '  +2.70%  '
$ws.Range('D3').Value = '2.054.98'
$ws.Range('E3').Value = '  +2.10%  '
$ws.Range('E4').Value = '  +0.09%  '
Set-TextValue 'D5' '230.03'
$ws.Range('E5').Value = '  +1.77%  '
$ws.Range('E6').Value = '  +1.64%  '
Set-TextValue 'D7' '58.35'
$ws.Range('E7').Value = '  +6.50%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('E9').Value = '  +2.86%  '
Set-TextValue 'D10' '0.0806'
$ws.Range('E10').Value = '  +3.21%  '
$ws.Range('E11').Value = '  +0.82%  '
$ws.Range('D12').Value = '2.360.30'
$ws.Range('E12').Value = '  +2.04%  '
Set-TextValue 'D13' '14.61'
$ws.Range('E13').Value = '  +3.58%  '
Set-TextValue 'D14' '20.66'
$ws.Range('E14').Value = '  +2.35%  '
Set-TextValue 'D15' '0.752'
$ws.Range('E15').Value = '  +1.57%  '
$ws.Range('E16').Value = '  +3.38%  '
$ws.Range('D17').Value = '2.060.66'
$ws.Range('E17').Value = '  +2.27%  '
$ws.Range('D18').Value = '37.939.59'
$ws.Range('E18').Value = '  +2.63%  '
Set-TextValue 'D19' '6.15'
$ws.Range('E19').Value = '  -0.32%  '
Set-TextValue 'D20' '69.82'
$ws.Range('E20').Value = '  +1.47%  '
$ws.Range('D21').Value = '0.0₃0829'
$ws.Range('E21').Value = '  +1.73%  '
Set-TextValue 'D22' '224.67'
$ws.Range('E22').Value = '  +0.87%  '
$ws.Range('E23').Value = '  +0.06%  '
$ws.Range('E24').Value = '  +1.11%  '
Set-TextValue 'D25' '2.25'
$ws.Range('E25').Value = '  +3.25%  '
$ws.Range('E26').Value = '  +1.74%  '
$ws.Range('E27').Value = '  +0.06%  '
Set-TextValue 'D28' '0.132'
$ws.Range('E28').Value = '  +6.67%  '
Set-TextValue 'D29' '19.02'
$ws.Range('E29').Value = '  +1.82%  '
$ws.Range('E30').Value = '  +0.13%  '
$ws.Range('E31').Value = '  +1.92%  '
$ws.Range('E32').Value = '  +0.60%  '
$ws.Range('E33').Value = '  +4.26%  '
Set-TextValue 'D34' '0.0612'
$ws.Range('E34').Value = '  +0.05%  '
$ws.Range('E35').Value = '  +8.06%  '
$ws.Range('E36').Value = '  +0.10%  '
Set-TextValue 'D37' '5.99'
$ws.Range('E37').Value = '  +13.22%  '
Set-TextValue 'D38' '3.32'
$ws.Range('E38').Value = '  +5.32%  '
$ws.Range('E39').Value = '  -0.22%  '
Set-TextValue 'D40' '98.43'
$ws.Range('E40').Value = '  +3.89%  '
$ws.Range('E41').Value = '  +1.80%  '
$ws.Range('D42').Value = '1.480.01'
$ws.Range('E42').Value = '  +0.12%  '
Set-TextValue 'D43' '0.0939'
$ws.Range('E43').Value = '  +2.89%  '
$ws.Range('E44').Value = '  +3.58%  '
Set-TextValue 'D45' '16.67'
$ws.Range('E45').Value = '  +2.61%  '
$ws.Range('E46').Value = '  +0.64%  '
Set-TextValue 'D47' '4.10'
$ws.Range('E47').Value = '  +16.54%  '
$ws.Range('E48').Value = '  +1.14%  '
$ws.Range('E49').Value = '  +1.62%  '
$ws.Range('E50').Value = '  -1.16%  '
$ws.Range('D51').Value = '2.249.12'
$ws.Range('E51').Value = '  +2.18%  '
